$wb = $excel.ActiveWorkbook

# --- Rename "Sub-tasks" sheet to "Subtasks" ---
# (This also auto-updates the _FilterDatabase defined name reference.)
$wsSub = $wb.Worksheets.Item("Sub-tasks")
$wsSub.Name = "Subtasks"

# --- Rename the "Sub-task" header label to "Subtask" ---
$wsSub.Range("A1").Value = "Subtask"

# --- Task-status sheet keeps its own selection at A2, but must be touched first so
#     that it ends up NOT being the active tab once Subtasks is activated below. ---
$wsStatus = $wb.Worksheets.Item("Task-status")
$wsStatus.Range("A2").Select()

# --- Reset selection on the Subtasks sheet to A2 and make it the active tab ---
$wsSub.Activate()
$wsSub.Range("A2").Select()
